# feat: add 2022-Q1 data
#
# Before: sheets = 2021-Q2, 2021-Q3, 2021-Q4, 总计
# After:  sheets = 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#   - "2022-Q1" is a brand-new fund-holdings sheet (same shape as the other
#     quarter sheets) inserted right before "总计".
#   - "总计" (the roll-up sheet) gains a new first data row for 2022-Q1 and
#     keeps its previous rows (shifted down).

$wb = $excel.ActiveWorkbook

# Chinese sheet name built from UTF-8 bytes so this script is encoding-safe
# regardless of how the host re-saves/reads this file.
$totalName = [System.Text.Encoding]::UTF8.GetString([byte[]](230,128,187,232,174,161))

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- Step 1: duplicate the "总计" sheet, placing the copy right after it.
#     Order becomes: ... 2021-Q4, 总计, 总计(2)
$total = $wb.Worksheets.Item($totalName)
$total.Copy([System.Reflection.Missing]::Value, $total)

# --- Step 2: duplicate "2021-Q4" (a fund-holding sheet, giving us the
#     right column layout/styles for the new quarter) placing it right
#     before the original "总计".
#     Order becomes: ... 2021-Q4, 2021-Q4(2), 总计, 总计(2)
$total = $wb.Worksheets.Item($totalName)
$q4.Copy($total, [System.Reflection.Missing]::Value)

# --- Step 3: drop the original "总计" sheet (we keep the later copy as the
#     real roll-up sheet so it ends up last).
$total = $wb.Worksheets.Item($totalName)
$excel.DisplayAlerts = $false
[void]$total.Delete()
$excel.DisplayAlerts = $true

# --- Step 4: rename the duplicated Q4 sheet to "2022-Q1" and update its
#     figures (fund codes/names/header stay identical to the other quarter
#     sheets - only the fund-size/position numbers differ).
$q1 = $wb.Worksheets.Item("2021-Q4" + " (2)")
$q1.Name = "2022-Q1"

$q1.Range("D2").Value = "'1.34"
$q1.Range("E2").Value = "'82.00"
$q1.Range("F2").Value = "'9.39"
$q1.Range("G2").Value = "'0.1258"

$q1.Range("D3").Value = "'1.34"
$q1.Range("E3").Value = "'82.00"
$q1.Range("F3").Value = "'9.39"
$q1.Range("G3").Value = "'0.1258"

# --- Step 5: rename the remaining "总计 (2)" copy back to "总计" and insert
#     the 2022-Q1 row at the top of the roll-up table.
$totalCopyName = $totalName + " (2)"
$finalTotal = $wb.Worksheets.Item($totalCopyName)
$finalTotal.Name = $totalName

$finalTotal = $wb.Worksheets.Item($totalName)
$finalTotal.Rows.Item(2).Insert()

# the inherited row-1 formatting on the blank inserted row should not carry
# over to the new data row's B:D cells (they are unstyled in the source data)
$finalTotal.Range("B2:D2").ClearFormats()

# give A2 the same style as the other index cells in column A (copy format
# only, so no value is disturbed)
$finalTotal.Range("A3").Copy()
$finalTotal.Range("A2").PasteSpecial(-4122)

$finalTotal.Range("A2").Value = 0
$finalTotal.Range("B2").Value = "2022-Q1"
$finalTotal.Range("C2").Value = 2
$finalTotal.Range("D2").Value = 0.25

# re-sequence the index column for the rows that shifted down
$finalTotal.Range("A3").Value = 1
$finalTotal.Range("A4").Value = 2
$finalTotal.Range("A5").Value = 3

Write-Host "Added 2022-Q1 sheet and updated totals."
